$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 (shifting nothing below it since it's the last row);
# this duplicates row 5's current formatting (including its date number
# format) into the freshly created row 6.
$ws.Rows("6:6").Insert(-4121)

# Row 5's date cell (A5) should now pick up the same number-format style
# that rows 3/4 use, so copy A3's format onto A5.
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Populate the new journal entry (Day 5 of N8N Learning).
$ws.Cells.Item(6, 1).Value = 45916
$ws.Cells.Item(6, 2).Value = "Citizen Complaint Response Automation"
$ws.Cells.Item(6, 3).Value = "Automates Typeform submissions by checking age eligibility and sending personalized Gmail responses based on complaint type."
$ws.Cells.Item(6, 4).Value = "Citizen Complaint Response Automation.json"

# Match the new selection recorded in the workbook.
$ws.Range("D12").Select()
